# feat (serializer): seralizer support Pair
#
# Inserts a new "Pair" array-of-pair column (V) right after the existing
# Map<uint> column (U), which itself is being re-purposed into a
# Pair<uint> column (tag format changes from "tag-N" to "tag:N").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a brand-new column at V; this pushes the old V..AG block
#    (the Array<float> nested-array columns) one column to the right,
#    becoming W..AH, and Excel keeps dimension/selection/col widths in
#    sync automatically.
$ws.Range("V1").EntireColumn.Insert()

# 2) Give the freshly inserted column a sensible width (matches the
#    widened "Array<Pair>" header text).
$ws.Range("V1").ColumnWidth = 18.33203125

# 3) Column U moves from "Map<uint>" semantics to "Pair<uint>" -- update
#    the type header, the type-name row, and the per-row tag labels.
$ws.Range("U5").Value = "Pair<uint>"
$ws.Range("U6").Value = "pair"

$ws.Range("U7").Value  = "tag:0"
$ws.Range("U8").Value  = "tag:1"
$ws.Range("U9").Value  = "tag:2"
$ws.Range("U10").Value = "tag:3"
$ws.Range("U11").Value = "tag:4"
$ws.Range("U12").Value = "tag:5"
$ws.Range("U13").Value = "tag:6"
$ws.Range("U14").Value = "tag:7"
$ws.Range("U15").Value = "tag:8"

# 4) Populate the new column V with the Array<Pair> data: header, type
#    name, and one "tag:sN" value per data row.
$ws.Range("V5").Value = "Array<Pair>"
$ws.Range("V6").Value = "map"

$ws.Range("V7").Value  = "tag:0"
$ws.Range("V8").Value  = "tag:s1"
$ws.Range("V9").Value  = "tag:s2"
$ws.Range("V10").Value = "tag:s3"
$ws.Range("V11").Value = "tag:s4"
$ws.Range("V12").Value = "tag:s5"
$ws.Range("V13").Value = "tag:s6"
$ws.Range("V14").Value = "tag:s7"
$ws.Range("V15").Value = "tag:s8"

# 5) Match the author's final selection.
$ws.Range("V10").Select()
